# Applies the authored diff to market_health_data.xlsx:
# - "Top Losers" sheet: rows re-ranked/updated (rows 18, 35-36, 48, 51-67, 70-76)
# - "1 Month Performance" sheet: rows re-ranked/updated (rows 6,9,13-14,29,36-40,52-62,75)
$wb = $excel.ActiveWorkbook

$wsLosers = $wb.Worksheets.Item("Top Losers")

$wsLosers.Range("D18").Value = -0.062
$wsLosers.Range("B35").Value = "SPARC"
$wsLosers.Range("C35").Value = -3.1709
$wsLosers.Range("D35").Value = 4.8337
$wsLosers.Range("E35").Value = 6.3311
$wsLosers.Range("B36").Value = "PRUDENT"
$wsLosers.Range("C36").Value = -3.127
$wsLosers.Range("D36").Value = -3.5103
$wsLosers.Range("E36").Value = 2.1213
$wsLosers.Range("D48").Value = 0.05
$wsLosers.Range("B51").Value = "UNIMECH"
$wsLosers.Range("C51").Value = -2.8008
$wsLosers.Range("D51").Value = -1.6104
$wsLosers.Range("E51").Value = -0.4585
$wsLosers.Range("B52").Value = "TTKPRESTIG"
$wsLosers.Range("C52").Value = -2.7438
$wsLosers.Range("D52").Value = 8.001200000000001
$wsLosers.Range("E52").Value = 9.650499999999999
$wsLosers.Range("B53").Value = "PFOCUS"
$wsLosers.Range("C53").Value = -2.7039
$wsLosers.Range("D53").Value = -2.6276
$wsLosers.Range("E53").Value = -1.2163
$wsLosers.Range("B54").Value = "ALLDIGI"
$wsLosers.Range("C54").Value = -2.6342
$wsLosers.Range("D54").Value = -0.2306
$wsLosers.Range("E54").Value = -5.3103
$wsLosers.Range("B55").Value = "PRIVISCL"
$wsLosers.Range("C55").Value = -2.6288
$wsLosers.Range("D55").Value = -2.1048
$wsLosers.Range("E55").Value = 19.7451
$wsLosers.Range("B56").Value = "CANHLIFE"
$wsLosers.Range("C56").Value = -2.6148
$wsLosers.Range("D56").Value = 3.7771
$wsLosers.Range("E56").Value = "N/A"
$wsLosers.Range("B57").Value = "GKENERGY"
$wsLosers.Range("C57").Value = -2.6122
$wsLosers.Range("D57").Value = -9.807700000000001
$wsLosers.Range("E57").Value = 23.2758
$wsLosers.Range("B58").Value = "SGFIN"
$wsLosers.Range("C58").Value = -2.592
$wsLosers.Range("D58").Value = -0.06270000000000001
$wsLosers.Range("E58").Value = 11.7235
$wsLosers.Range("B59").Value = "ARVINDFASN"
$wsLosers.Range("C59").Value = -2.549
$wsLosers.Range("D59").Value = -2.9892
$wsLosers.Range("E59").Value = -4.4223
$wsLosers.Range("B60").Value = "EDELWEISS"
$wsLosers.Range("C60").Value = -2.5422
$wsLosers.Range("D60").Value = -3.3745
$wsLosers.Range("E60").Value = 8.5305
$wsLosers.Range("B61").Value = "SAMHI"
$wsLosers.Range("C61").Value = -2.5284
$wsLosers.Range("D61").Value = 1.8231
$wsLosers.Range("E61").Value = 2.8516
$wsLosers.Range("B62").Value = "TBOTEK"
$wsLosers.Range("C62").Value = -2.524
$wsLosers.Range("D62").Value = -3.5732
$wsLosers.Range("E62").Value = 1.036
$wsLosers.Range("B63").Value = "UJJIVANSFB"
$wsLosers.Range("C63").Value = -2.5201
$wsLosers.Range("D63").Value = 0.3845
$wsLosers.Range("E63").Value = 12.6645
$wsLosers.Range("B64").Value = "AMBER"
$wsLosers.Range("C64").Value = -2.5098
$wsLosers.Range("D64").Value = -0.1082
$wsLosers.Range("E64").Value = 2.763
$wsLosers.Range("B65").Value = "GRPLTD"
$wsLosers.Range("C65").Value = -2.4898
$wsLosers.Range("D65").Value = -5.9894
$wsLosers.Range("E65").Value = -5.4586
$wsLosers.Range("B66").Value = "NESCO"
$wsLosers.Range("C66").Value = -2.4722
$wsLosers.Range("D66").Value = 1.9934
$wsLosers.Range("E66").Value = 3.8931
$wsLosers.Range("B67").Value = "PILANIINVS"
$wsLosers.Range("C67").Value = -2.4546
$wsLosers.Range("D67").Value = -0.7907
$wsLosers.Range("E67").Value = 4.267
$wsLosers.Range("B70").Value = "JNKINDIA"
$wsLosers.Range("C70").Value = -2.3482
$wsLosers.Range("D70").Value = -2.8371
$wsLosers.Range("E70").Value = 4.2622
$wsLosers.Range("B71").Value = "FCL"
$wsLosers.Range("C71").Value = -2.3453
$wsLosers.Range("D71").Value = -2.616
$wsLosers.Range("E71").Value = -0.02
$wsLosers.Range("B72").Value = "DEEDEV"
$wsLosers.Range("C72").Value = -2.3334
$wsLosers.Range("D72").Value = -6.6528
$wsLosers.Range("E72").Value = -7.4227
$wsLosers.Range("B73").Value = "WEALTH"
$wsLosers.Range("C73").Value = -2.2793
$wsLosers.Range("D73").Value = -3.8356
$wsLosers.Range("E73").Value = -2.7981
$wsLosers.Range("B74").Value = "RATNAMANI"
$wsLosers.Range("C74").Value = -2.2788
$wsLosers.Range("D74").Value = -0.4626
$wsLosers.Range("E74").Value = 0.8712
$wsLosers.Range("B75").Value = "CSBBANK"
$wsLosers.Range("C75").Value = -2.2695
$wsLosers.Range("D75").Value = 2.3137
$wsLosers.Range("E75").Value = 10.6999
$wsLosers.Range("B76").Value = "BBOX"
$wsLosers.Range("C76").Value = -2.2639
$wsLosers.Range("D76").Value = -4.7636
$wsLosers.Range("E76").Value = 5.1528

$wsPerf = $wb.Worksheets.Item("1 Month Performance")

$wsPerf.Range("C6").Value = 66.1892
$wsPerf.Range("C9").Value = 55.9703
$wsPerf.Range("B13").Value = "TVSSRICHAK"
$wsPerf.Range("C13").Value = 40.7778
$wsPerf.Range("B14").Value = "MTARTECH"
$wsPerf.Range("C14").Value = 40.7213
$wsPerf.Range("C29").Value = 31.8611
$wsPerf.Range("B36").Value = "MINDTECK"
$wsPerf.Range("C36").Value = 26.9415
$wsPerf.Range("B37").Value = "BHARATWIRE"
$wsPerf.Range("C37").Value = 26.5276
$wsPerf.Range("B38").Value = "HATSUN"
$wsPerf.Range("C38").Value = 26.492
$wsPerf.Range("B39").Value = "INDORAMA"
$wsPerf.Range("C39").Value = 26.4516
$wsPerf.Range("B40").Value = "IFBIND"
$wsPerf.Range("C40").Value = 26.161
$wsPerf.Range("B52").Value = "PVP"
$wsPerf.Range("C52").Value = 22.7524
$wsPerf.Range("B53").Value = "INDIANB"
$wsPerf.Range("C53").Value = 22.6689
$wsPerf.Range("B54").Value = "PRIVISCL"
$wsPerf.Range("C54").Value = 22.3984
$wsPerf.Range("B55").Value = "CPEDU"
$wsPerf.Range("C55").Value = 22.3786
$wsPerf.Range("B56").Value = "LORDSCHLO"
$wsPerf.Range("C56").Value = 22.1791
$wsPerf.Range("B57").Value = "GUJTHEM"
$wsPerf.Range("C57").Value = 22.0704
$wsPerf.Range("B58").Value = "SURYODAY"
$wsPerf.Range("C58").Value = 21.8039
$wsPerf.Range("B59").Value = "TDPOWERSYS"
$wsPerf.Range("C59").Value = 21.7743
$wsPerf.Range("B60").Value = "ORBTEXP"
$wsPerf.Range("C60").Value = 21.6115
$wsPerf.Range("B61").Value = "CEATLTD"
$wsPerf.Range("C61").Value = 20.0239
$wsPerf.Range("B62").Value = "ATL"
$wsPerf.Range("C62").Value = 19.9362
$wsPerf.Range("B75").Value = "HITECHGEAR"
$wsPerf.Range("C75").Value = 18.6598
